$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readme")
$c = $ws.Cells.Item(2,5)
$c.Value(1) = "20220302"
